$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H3").Value = 24551.5
$ws.Range("J3").Value = 24551.5
$ws.Range("L3").Value = 24551.5
$ws.Range("N3").Value = -24779.5

$ws = $wb.Worksheets.Item(1)
$ws.Range("H11").Value = 10.1875
$ws.Range("I11").Value = 10.1875
$ws.Range("K11").Value = 10.1875
$ws.Range("M11").Value = 129.8125

$ws = $wb.Worksheets.Item(1)
$ws.Range("H41").Value = 2332
$ws.Range("I41").Value = 2332
$ws.Range("K41").Value = 2332
$ws.Range("M41").Value = -1892

$ws = $wb.Worksheets.Item(1)
$ws.Range("H58").Value = 1896.0526
$ws.Range("J58").Value = 2386.3333
$ws.Range("L58").Value = 7158.999899999999
$ws.Range("N58").Value = -7458.999899999999

$ws = $wb.Worksheets.Item(1)
$ws.Range("H102").Value = 24551.5
$ws.Range("J102").Value = 24551.5
$ws.Range("L102").Value = 24551.5
$ws.Range("N102").Value = -31041.5

$ws = $wb.Worksheets.Item(1)
$ws.Range("H111").Value = 1085.25
$ws.Range("J111").Value = 1270.5
$ws.Range("L111").Value = 3811.5
$ws.Range("N111").Value = -9945.5

$ws = $wb.Worksheets.Item(1)
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 2755.818
$ws.Range("I74").Value = 2221.9
$ws.Range("K74").Value = 2221.9
$ws.Range("M74").Value = -1347.9

$ws = $wb.Worksheets.Item(2)
$ws.Range("H77").Value = 2755.818
$ws.Range("I77").Value = 2221.9
$ws.Range("K77").Value = 11109.5
$ws.Range("M77").Value = -6741.5

$ws = $wb.Worksheets.Item(2)
$ws.Range("H102").Value = 6254059.5
$ws.Range("I102").Value = 8930971
$ws.Range("J102").Value = 7931.6665
$ws.Range("K102").Value = 8930971
$ws.Range("L102").Value = 7931.6665
$ws.Range("M102").Value = -8929349
$ws.Range("N102").Value = -11175.6665

$ws = $wb.Worksheets.Item(2)
$ws.Range("H110").Value = 111114456
$ws.Range("J110").Value = 4015.25
$ws.Range("L110").Value = 4015.25
$ws.Range("N110").Value = -8105.25

$ws = $wb.Worksheets.Item(2)
$ws.Range("H122").Value = 2553.5833
$ws.Range("I122").Value = 1288
$ws.Range("J122").Value = 4325.4
$ws.Range("K122").Value = 3864
$ws.Range("L122").Value = 12976.2
$ws.Range("M122").Value = -1414
$ws.Range("N122").Value = -17876.2

$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 3005.6667
$ws.Range("I20").Value = 4004
$ws.Range("K20").Value = 4004
$ws.Range("M20").Value = -3757

$ws = $wb.Worksheets.Item(3)
$ws.Range("H63").Value = 74749.5
$ws.Range("J63").Value = 66333
$ws.Range("L63").Value = 66333
$ws.Range("N63").Value = -67705

$ws = $wb.Worksheets.Item(3)
$ws.Range("H66").Value = 74749.5
$ws.Range("J66").Value = 66333
$ws.Range("L66").Value = 198999
$ws.Range("N66").Value = -205863

$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 6342.75
$ws.Range("I86").Value = 3651.5
$ws.Range("J86").Value = 7688.375
$ws.Range("K86").Value = 3651.5
$ws.Range("L86").Value = 7688.375
$ws.Range("M86").Value = -2528.5
$ws.Range("N86").Value = -9934.375

$ws = $wb.Worksheets.Item(3)
$ws.Range("H89").Value = 6342.75
$ws.Range("I89").Value = 3651.5
$ws.Range("J89").Value = 7688.375
$ws.Range("K89").Value = 18257.5
$ws.Range("L89").Value = 38441.875
$ws.Range("M89").Value = -12641.5
$ws.Range("N89").Value = -49673.875

$ws = $wb.Worksheets.Item(3)
$ws.Range("H94").Value = 493.83334
$ws.Range("I94").Value = 482.8
$ws.Range("J94").Value = 549
$ws.Range("K94").Value = 482.8
$ws.Range("L94").Value = 549
$ws.Range("M94").Value = -31.80000000000001
$ws.Range("N94").Value = -1451

$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 3900.8572
$ws.Range("I134").Value = 1043.0526
$ws.Range("J134").Value = 31050
$ws.Range("K134").Value = 3129.1578
$ws.Range("L134").Value = 93150
$ws.Range("M134").Value = -594.1578
$ws.Range("N134").Value = -98220

$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 5838.069
$ws.Range("I31").Value = 2925.25
$ws.Range("K31").Value = 2925.25
$ws.Range("M31").Value = -2630.25

$ws = $wb.Worksheets.Item(4)
$ws.Range("H34").Value = 5838.069
$ws.Range("I34").Value = 2925.25
$ws.Range("K34").Value = 2925.25
$ws.Range("M34").Value = -2723.25

$ws = $wb.Worksheets.Item(4)
$ws.Range("H45").Value = 24999.5
$ws.Range("I45").Value = 10000
$ws.Range("J45").Value = 39999
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 39999
$ws.Range("M45").Value = -9407
$ws.Range("N45").Value = -41185

$ws = $wb.Worksheets.Item(4)
$ws.Range("H106").Value = 38199.5
$ws.Range("J106").Value = 38199.5
$ws.Range("L106").Value = 38199.5
$ws.Range("N106").Value = -40723.5

$ws = $wb.Worksheets.Item(4)
$ws.Range("H132").Value = 1415
$ws.Range("I132").Value = 1415
$ws.Range("K132").Value = 4245
$ws.Range("M132").Value = -1715

$ws = $wb.Worksheets.Item(5)
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws = $wb.Worksheets.Item(5)
$ws.Range("H26").Value = 168.25
$ws.Range("J26").Value = 100
$ws.Range("L26").Value = 300
$ws.Range("N26").Value = -876

$ws = $wb.Worksheets.Item(5)
$ws.Range("H80").Value = 4235.2593
$ws.Range("I80").Value = 4020.9473
$ws.Range("J80").Value = 4744.25
$ws.Range("K80").Value = 12062.8419
$ws.Range("L80").Value = 14232.75
$ws.Range("M80").Value = -11126.8419
$ws.Range("N80").Value = -16104.75

$ws = $wb.Worksheets.Item(5)
$ws.Range("H83").Value = 4235.2593
$ws.Range("I83").Value = 4020.9473
$ws.Range("J83").Value = 4744.25
$ws.Range("K83").Value = 36188.5257
$ws.Range("L83").Value = 42698.25
$ws.Range("M83").Value = -31508.5257
$ws.Range("N83").Value = -52058.25

$ws = $wb.Worksheets.Item(5)
$ws.Range("H132").Value = 3220.75
$ws.Range("I132").Value = 1900
$ws.Range("J132").Value = 3661
$ws.Range("K132").Value = 17100
$ws.Range("L132").Value = 32949
$ws.Range("M132").Value = -14570
$ws.Range("N132").Value = -38009

$ws = $wb.Worksheets.Item(5)
$ws.Range("H137").Value = 5761.25
$ws.Range("J137").Value = 5761.25
$ws.Range("L137").Value = 17283.75
$ws.Range("N137").Value = -27483.75

$ws = $wb.Worksheets.Item(6)
$ws.Range("H49").Value = 23997
$ws.Range("J49").Value = 23997
$ws.Range("L49").Value = 23997
$ws.Range("N49").Value = -24365

$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws = $wb.Worksheets.Item(6)
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 52055.55
$ws.Range("I132").Value = 52055.55
$ws.Range("K132").Value = 156166.65
$ws.Range("M132").Value = -153636.65

$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 2917.4
$ws.Range("I16").Value = 2917.4
$ws.Range("K16").Value = 2917.4
$ws.Range("M16").Value = -2747.4

$ws = $wb.Worksheets.Item(7)
$ws.Range("H93").Value = 1439.375
$ws.Range("I93").Value = 1379.2307
$ws.Range("J93").Value = 1700
$ws.Range("K93").Value = 1379.2307
$ws.Range("L93").Value = 1700
$ws.Range("M93").Value = -131.2307000000001
$ws.Range("N93").Value = -4196

$ws = $wb.Worksheets.Item(7)
$ws.Range("H122").Value = 4995.6
$ws.Range("I122").Value = 4996
$ws.Range("K122").Value = 14988
$ws.Range("M122").Value = -12538

$ws = $wb.Worksheets.Item(7)
$ws.Range("H132").Value = 3909.2
$ws.Range("I132").Value = 3886.5
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 11659.5
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -9129.5
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item(8)
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 1138.4286
$ws.Range("I132").Value = 1228.1666
$ws.Range("J132").Value = 600
$ws.Range("K132").Value = 3684.4998
$ws.Range("L132").Value = 1800
$ws.Range("M132").Value = -1154.4998
$ws.Range("N132").Value = -6860

$ws = $wb.Worksheets.Item(8)
$ws.Range("H136").Value = 2716.423
$ws.Range("I136").Value = 1786.3125
$ws.Range("K136").Value = 5358.9375
$ws.Range("M136").Value = -2808.9375
